$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts old rows 3-10 down to 4-11),
# matching the formatting of the row above it.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with this week's entry (copy of the former
# row 3 data, but with the new date and a relabeled commercialization unit).
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44545
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 23500
$ws.Range("Q3").Value = "`$/bandeja 12 kilos"
$ws.Range("R3").Value = "Región de Coquimbo"
$ws.Range("S3").Value = 1958
$ws.Range("T3").Value = 12
